# Auto-generated update of Leve price/profit data (scheduled runner refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3015
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 3015
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 9045
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -10793
$ws.Range("H72").Value = 3015
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 3015
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 27135
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -35871
$ws.Range("H87").Value = 29838
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 29838
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 29838
$ws.Range("N87").Value = -32334
$ws.Range("H90").Value = 29838
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 29838
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 89514
$ws.Range("N90").Value = -101994
$ws.Range("H114").Value = 36786
$ws.Range("I114").Value = 36786
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 36786
$ws.Range("L114").ClearContents()
$ws.Range("M114").Value = -32447
$ws.Range("N114").Value = 0
$ws.Range("H129").Value = 970.0137
$ws.Range("I129").Value = 974.9
$ws.Range("J129").Value = 969.2381
$ws.Range("K129").Value = 2924.7
$ws.Range("L129").Value = 2907.7143
$ws.Range("M129").Value = 2075.3
$ws.Range("N129").Value = -12907.7143
$ws.Range("H138").Value = 3088864
$ws.Range("I138").Value = 6500
$ws.Range("J138").Value = 3207416.5
$ws.Range("K138").Value = 19500
$ws.Range("L138").Value = 9622249.5
$ws.Range("M138").Value = -14360
$ws.Range("N138").Value = -9632529.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 143143650
$ws.Range("I61").Value = 200200660
$ws.Range("J61").Value = 501111
$ws.Range("K61").Value = 200200660
$ws.Range("L61").Value = 501111
$ws.Range("M61").Value = -200200448
$ws.Range("N61").Value = -501535
$ws.Range("H122").Value = 1902
$ws.Range("I122").Value = 1902
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5706
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -3256
$ws.Range("H136").Value = 143143650
$ws.Range("I136").Value = 200200660
$ws.Range("J136").Value = 501111
$ws.Range("K136").Value = 600601980
$ws.Range("L136").Value = 1503333
$ws.Range("M136").Value = -600599430
$ws.Range("N136").Value = -1508433

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2049.2856
$ws.Range("I134").Value = 2025.2632
$ws.Range("J134").Value = 2100
$ws.Range("K134").Value = 6075.7896
$ws.Range("L134").Value = 6300
$ws.Range("M134").Value = -3540.7896
$ws.Range("N134").Value = -11370

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1753.5883
$ws.Range("I122").Value = 1412.7368
$ws.Range("J122").Value = 2185.3333
$ws.Range("K122").Value = 4238.2104
$ws.Range("L122").Value = 6555.999899999999
$ws.Range("M122").Value = -1788.2104
$ws.Range("N122").Value = -11455.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4451.2856
$ws.Range("I70").Value = 4436.8887
$ws.Range("J70").Value = 4477.2
$ws.Range("K70").Value = 4436.8887
$ws.Range("L70").Value = 4477.2
$ws.Range("M70").Value = -4166.8887
$ws.Range("N70").Value = -5017.2
$ws.Range("H73").Value = 4451.2856
$ws.Range("I73").Value = 4436.8887
$ws.Range("J73").Value = 4477.2
$ws.Range("K73").Value = 4436.8887
$ws.Range("L73").Value = 4477.2
$ws.Range("M73").Value = -3500.8887
$ws.Range("N73").Value = -6349.2
$ws.Range("H80").Value = 3043.3333
$ws.Range("I80").Value = 2281.1538
$ws.Range("J80").Value = 3751.0715
$ws.Range("K80").Value = 2281.1538
$ws.Range("L80").Value = 3751.0715
$ws.Range("M80").Value = -1283.1538
$ws.Range("N80").Value = -5747.0715
$ws.Range("H83").Value = 3043.3333
$ws.Range("I83").Value = 2281.1538
$ws.Range("J83").Value = 3751.0715
$ws.Range("K83").Value = 11405.769
$ws.Range("L83").Value = 18755.3575
$ws.Range("M83").Value = -6413.769
$ws.Range("N83").Value = -28739.3575
$ws.Range("H122").Value = 1611.5714
$ws.Range("I122").Value = 1463.5
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 4390.5
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -1940.5
$ws.Range("N122").Value = -12400
$ws.Range("H132").Value = 253188.75
$ws.Range("I132").Value = 334004
$ws.Range("J132").Value = 204699.6
$ws.Range("K132").Value = 1002012
$ws.Range("L132").Value = 614098.8
$ws.Range("M132").Value = -999482
$ws.Range("N132").Value = -619158.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 708.25
$ws.Range("I22").Value = 415.5
$ws.Range("J22").Value = 1001
$ws.Range("K22").Value = 415.5
$ws.Range("L22").Value = 1001
$ws.Range("M22").Value = -120.5
$ws.Range("N22").Value = -1591
$ws.Range("H27").Value = 708.25
$ws.Range("I27").Value = 415.5
$ws.Range("J27").Value = 1001
$ws.Range("K27").Value = 415.5
$ws.Range("L27").Value = 1001
$ws.Range("M27").Value = -308.5
$ws.Range("N27").Value = -1215
$ws.Range("H68").Value = 1622.35
$ws.Range("I68").Value = 1574.7858
$ws.Range("J68").Value = 1733.3334
$ws.Range("K68").Value = 1574.7858
$ws.Range("L68").Value = 1733.3334
$ws.Range("M68").Value = -825.7858000000001
$ws.Range("N68").Value = -3231.3334
$ws.Range("H71").Value = 1622.35
$ws.Range("I71").Value = 1574.7858
$ws.Range("J71").Value = 1733.3334
$ws.Range("K71").Value = 7873.929
$ws.Range("L71").Value = 8666.666999999999
$ws.Range("M71").Value = -4129.929
$ws.Range("N71").Value = -16154.667
$ws.Range("H122").Value = 4080.6428
$ws.Range("I122").Value = 3397.3333
$ws.Range("J122").Value = 4593.125
$ws.Range("K122").Value = 10191.9999
$ws.Range("L122").Value = 13779.375
$ws.Range("M122").Value = -7741.999899999999
$ws.Range("N122").Value = -18679.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8001.5
$ws.Range("I62").Value = 6000
$ws.Range("J62").Value = 10003
$ws.Range("K62").Value = 6000
$ws.Range("L62").Value = 10003
$ws.Range("M62").Value = -5376
$ws.Range("N62").Value = -11251
$ws.Range("H65").Value = 8001.5
$ws.Range("I65").Value = 6000
$ws.Range("J65").Value = 10003
$ws.Range("K65").Value = 30000
$ws.Range("L65").Value = 50015
$ws.Range("M65").Value = -26880
$ws.Range("N65").Value = -56255
$ws.Range("H81").Value = 2583.1667
$ws.Range("I81").Value = 2000
$ws.Range("J81").Value = 2874.75
$ws.Range("K81").Value = 4000
$ws.Range("L81").Value = 5749.5
$ws.Range("M81").Value = -2939
$ws.Range("N81").Value = -7871.5
$ws.Range("H84").Value = 2583.1667
$ws.Range("I84").Value = 2000
$ws.Range("J84").Value = 2874.75
$ws.Range("K84").Value = 20000
$ws.Range("L84").Value = 28747.5
$ws.Range("M84").Value = -14696
$ws.Range("N84").Value = -39355.5
$ws.Range("H122").Value = 2443.158
$ws.Range("I122").Value = 2077.0833
$ws.Range("J122").Value = 3070.7144
$ws.Range("K122").Value = 6231.249899999999
$ws.Range("L122").Value = 9212.143199999999
$ws.Range("M122").Value = -3781.249899999999
$ws.Range("N122").Value = -14112.1432
